$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.474.58"
$ws.Range("E2").Value = "  -1.07%  "
$ws.Range("D3").Value = "2.062.45"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.90"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.41"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.17%  "
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0775"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("D12").Value = "2.366.98"
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.70"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.07"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.762"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.35"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "2.076.29"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "37.470.13"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.13"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.79"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("D21").Value = "0.0₃0827"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.21"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.94"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +8.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.64"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.131"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.25"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("E30").Value = "  -4.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.122"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.55"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0624"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.30"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("E40").Value = "  +3.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.82"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("E42").Value = "  +4.87%  "
$ws.Range("E43").Value = "  -2.37%  "
$ws.Range("D44").Value = "1.478.87"
$ws.Range("E44").Value = "  +2.35%  "
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.80"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("E47").Value = "  -2.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.98"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.65%  "
$ws.Range("E49").Value = "  -1.78%  "
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").Value = "2.251.37"
$ws.Range("E51").Value = "  -1.10%  "
